# Applies the diff: a new data row is inserted above former row 271,
# pushing the existing rows 271:313 down to 272:314 (keeping all their
# original values), and the newly inserted row 271 receives the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 271 (shifts rows 271:313 down to 272:314).
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new record's data.
$ws.Cells.Item(271, 1).Value = 4
$ws.Cells.Item(271, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(271, 3).Value = "Los Lagos"
$ws.Cells.Item(271, 4).Value = 44694
$ws.Cells.Item(271, 5).Value = 10
$ws.Cells.Item(271, 6).Value = 100112008
$ws.Cells.Item(271, 7).Value = "Coliflor"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Segunda"
$ws.Cells.Item(271, 10).Value = 1000
$ws.Cells.Item(271, 11).Value = 1500
$ws.Cells.Item(271, 12).Value = 1500
$ws.Cells.Item(271, 13).Value = 1500
$ws.Cells.Item(271, 14).Value = "`$/unidad"
$ws.Cells.Item(271, 15).Value = "Región Metropolitana"
$ws.Cells.Item(271, 16).Value = 1500
$ws.Cells.Item(271, 17).Value = 1
$ws.Cells.Item(271, 18).Value = "Hortaliza"

# Match the date cell format used by the other rows in column D.
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(272, 4).NumberFormat
